# Apply narrow (0.5") margins to every section and widen the document's
# tables so their columns continue to span the full (now-wider) content
# area between the margins.
#
# 0.5 in = 36 pt = 720 twips  (Word COM page-setup properties are in points)
$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $sec.PageSetup.TopMargin    = 36
    $sec.PageSetup.BottomMargin = 36
    $sec.PageSetup.LeftMargin   = 36
    $sec.PageSetup.RightMargin  = 36
}

# New content width = 12240 (8.5") page - 2*720 (0.5" margins) = 10800 twips
# = 270 pt per column for the two 2-column tables in this document.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    for ($j = 1; $j -le $tbl.Columns.Count; $j++) {
        $tbl.Columns.Item($j).Width = 270
    }
}
